$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original styling of the data range, force text format so
# numeric-looking strings (prices like "29.294.66") are not coerced to numbers,
# then restore the original style afterwards so appearance/format is unchanged.
$dataRange = $ws.Range("B2:E51")
$origStyle = $dataRange.Style
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "29.294.66"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "1.860.16"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("D4").Value = "0.9977"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "0.7114"
$ws.Range("E5").Value = "  +2.76%  "
$ws.Range("D6").Value = "237.92"
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("D7").Value = "0.9987"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").Value = "0.07964"
$ws.Range("E8").Value = "  +3.04%  "
$ws.Range("D9").Value = "0.3041"
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").Value = "23.58"
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("E11").Value = "  +0.38%  "
$ws.Range("D12").Value = "1.845.36"
$ws.Range("E12").Value = "  -0.93%  "
$ws.Range("D13").Value = "5.179"
$ws.Range("E13").Value = "  -1.09%  "
$ws.Range("D14").Value = "0.7048"
$ws.Range("E14").Value = "  -3.08%  "
$ws.Range("D15").Value = "89.79"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").Value = "29.233.68"
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("D17").Value = "5.853"
$ws.Range("E17").Value = "  +0.86%  "
$ws.Range("D18").Value = "0.000007896"
$ws.Range("E18").Value = "  +1.43%  "
$ws.Range("D19").Value = "13.31"
$ws.Range("E19").Value = "  +0.85%  "
$ws.Range("D20").Value = "237.82"
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("D21").Value = "0.9982"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").Value = "2.076.20"
$ws.Range("E22").Value = "  -1.69%  "
$ws.Range("D23").Value = "0.9977"
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("D24").Value = "7.412"
$ws.Range("E24").Value = "  -2.96%  "
$ws.Range("D25").Value = "162.10"
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("D26").Value = "8.946"
$ws.Range("E26").Value = "  -0.82%  "
$ws.Range("D27").Value = "0.1437"
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("D28").Value = "18.07"
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").Value = "1.929"
$ws.Range("E29").Value = "  -2.69%  "
$ws.Range("D30").Value = "1.433"
$ws.Range("E30").Value = "  +2.39%  "
$ws.Range("D31").Value = "1.482"
$ws.Range("D32").Value = "4.380"
$ws.Range("E32").Value = "  -2.29%  "
$ws.Range("D33").Value = "4.015"
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("D34").Value = "0.05216"
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("D35").Value = "1.164"
$ws.Range("E35").Value = "  -2.21%  "
$ws.Range("D36").Value = "0.7065"
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("D37").Value = "1.002"
$ws.Range("E37").Value = "  -3.40%  "
$ws.Range("D38").Value = "2.662"
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("D39").Value = "0.01850"
$ws.Range("E39").Value = "  -0.64%  "
$ws.Range("D40").Value = "2.723"
$ws.Range("E40").Value = "  +1.74%  "
$ws.Range("D41").Value = "0.9294"
$ws.Range("E41").Value = "  -0.85%  "
$ws.Range("D42").Value = "1.125.45"
$ws.Range("E42").Value = "  +4.15%  "
$ws.Range("D43").Value = "0.4268"
$ws.Range("E43").Value = "  -0.39%  "
$ws.Range("D44").Value = "5.845"
$ws.Range("E44").Value = "  -3.20%  "
$ws.Range("D45").Value = "69.83"
$ws.Range("E45").Value = "  -0.90%  "
$ws.Range("D46").Value = "0.9984"
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("D47").Value = "103.10"
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "0.5339"
$ws.Range("E48").Value = "  -4.29%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "1.763"
$ws.Range("E49").Value = "  -1.26%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "1.989.75"
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "9.131"
$ws.Range("E51").Value = "  -0.98%  "

# Restore original style/number format on the whole data range
$dataRange.Style = $origStyle
